$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO" (columns A:R) - insert a new advisor row
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row above the current row 2 ("ANGULO PARRALES CARMEN")
$ws1.Range("A2").EntireRow.Insert()

# Copy the formatting from the row that just got pushed down to row 3
# onto the newly inserted (blank) row 2, so number formats / styles match
# the other data rows instead of inheriting the header row's style.
$ws1.Range("A3:R3").Copy()
$ws1.Range("A2:R2").PasteSpecial(-4122)

# Fill in the new advisor's data
$ws1.Range("A2").Value = "OFICINA-CATAECSA"
$ws1.Range("B2").Value = "ALCIVAR BUSTAMANTE ERNESTO EDUARDO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(2, $col).Value = 0
}

# Update the trailing "X de 8" summary row (now row 11) to "X de 9"
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(11, $col)
    $cell.Value = ($cell.Value2 -replace "de 8", "de 9")
}

# Widen the CLIENTE column
$ws1.Columns.Item(2).ColumnWidth = 35.17

# ----------------------------------------------------------------------
# Sheet "VENTA MENSUAL" (columns A:G) - insert the same new advisor row
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("A2").EntireRow.Insert()

$ws2.Range("A3:G3").Copy()
$ws2.Range("A2:G2").PasteSpecial(-4122)

$ws2.Range("A2").Value = "OFICINA-CATAECSA"
$ws2.Range("B2").Value = "ALCIVAR BUSTAMANTE ERNESTO EDUARDO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(2, $col).Value = 0
}

$ws2.Columns.Item(2).ColumnWidth = 35.17
